$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.192379144208543
$ws.Cells.Item(2, 3).Value = 1.020878496217309
$ws.Cells.Item(3, 2).Value = 1.869410725055877
$ws.Cells.Item(3, 3).Value = -1.281648130147154
$ws.Cells.Item(4, 2).Value = -0.3318649995207109
$ws.Cells.Item(4, 3).Value = 0.5050772145142584
$ws.Cells.Item(5, 2).Value = 0.5404290999940266
$ws.Cells.Item(5, 3).Value = 0.3944810077655069
$ws.Cells.Item(6, 2).Value = -2.075031836665305
$ws.Cells.Item(6, 3).Value = 0.6967747618062374
$ws.Cells.Item(7, 2).Value = -0.5111026748439523
$ws.Cells.Item(7, 3).Value = 1.258140074602792
$ws.Cells.Item(8, 2).Value = 0.3199689044927579
$ws.Cells.Item(8, 3).Value = 0.6547783786586647
$ws.Cells.Item(9, 2).Value = -1.35617258160237
$ws.Cells.Item(9, 3).Value = 0.2897668936346211
$ws.Cells.Item(10, 2).Value = -1.759167151189836
$ws.Cells.Item(10, 3).Value = 0.218972355154201
$ws.Cells.Item(11, 2).Value = 0.2081726550663906
$ws.Cells.Item(11, 3).Value = -0.09762189235145383
$ws.Cells.Item(12, 2).Value = 0.5739423672284978
$ws.Cells.Item(12, 3).Value = 0.1173796061076274
$ws.Cells.Item(13, 2).Value = 0.5783484400170088
$ws.Cells.Item(13, 3).Value = 2.253570470933922
$ws.Cells.Item(14, 2).Value = -1.056023949464658
$ws.Cells.Item(14, 3).Value = 1.521058393174677
$ws.Cells.Item(15, 2).Value = 1.535231403256816
$ws.Cells.Item(15, 3).Value = -0.1210661790642925
$ws.Cells.Item(16, 2).Value = 1.386474917277614
$ws.Cells.Item(16, 3).Value = 0.4228447078699191
$ws.Cells.Item(17, 2).Value = 1.622349763539734
$ws.Cells.Item(17, 3).Value = -1.655081222885729
$ws.Cells.Item(18, 2).Value = 0.07375551398390523
$ws.Cells.Item(18, 3).Value = 2.744087149141673
$ws.Cells.Item(19, 2).Value = 0.5887235888409346
$ws.Cells.Item(19, 3).Value = -1.503518955136582
$ws.Cells.Item(20, 2).Value = 0.2329519635282597
$ws.Cells.Item(20, 3).Value = 1.312606378596723
$ws.Cells.Item(21, 2).Value = 0.6977626795966044
$ws.Cells.Item(21, 3).Value = -0.4454913566412553
$ws.Cells.Item(22, 2).Value = -0.00007823961179193332
$ws.Cells.Item(22, 3).Value = 0.8118579174230005
$ws.Cells.Item(23, 2).Value = -0.3209752567494464
$ws.Cells.Item(23, 3).Value = 1.32721658784677
$ws.Cells.Item(24, 2).Value = -2.731244840337007
$ws.Cells.Item(24, 3).Value = 0.2206305342992531
$ws.Cells.Item(25, 2).Value = -0.6388527528244204
$ws.Cells.Item(25, 3).Value = 0.6918707534773548
$ws.Cells.Item(26, 2).Value = -0.2368370806480016
$ws.Cells.Item(26, 3).Value = -1.098004355100295
$ws.Cells.Item(27, 2).Value = 1.340773809355093
$ws.Cells.Item(27, 3).Value = -0.5214057232960154
$ws.Cells.Item(28, 2).Value = 0.1359954553019375
$ws.Cells.Item(28, 3).Value = -1.475737748515807
$ws.Cells.Item(29, 2).Value = -0.991772331957005
$ws.Cells.Item(29, 3).Value = -0.3836665694253608
$ws.Cells.Item(30, 2).Value = 1.611006750604784
$ws.Cells.Item(30, 3).Value = 0.2387426199344161
$ws.Cells.Item(31, 2).Value = 1.13992759312113
$ws.Cells.Item(31, 3).Value = -1.050319633043887
$ws.Cells.Item(32, 2).Value = -1.833060252516225
$ws.Cells.Item(32, 3).Value = -1.08054182651166
$ws.Cells.Item(33, 2).Value = 0.3698471364797039
$ws.Cells.Item(33, 3).Value = -0.5544481275710363
$ws.Cells.Item(34, 2).Value = 0.4843941365775502
$ws.Cells.Item(34, 3).Value = 0.4540101986283487
$ws.Cells.Item(35, 2).Value = 1.312609627433891
$ws.Cells.Item(35, 3).Value = -0.9707706091423735
$ws.Cells.Item(36, 2).Value = 0.2643944700618042
$ws.Cells.Item(36, 3).Value = 1.052694640878439
$ws.Cells.Item(37, 2).Value = 0.917289058043892
$ws.Cells.Item(37, 3).Value = -0.02152797604127773
$ws.Cells.Item(38, 2).Value = -1.217913137748882
$ws.Cells.Item(38, 3).Value = -1.74789066163599
$ws.Cells.Item(39, 2).Value = -1.248966940002415
$ws.Cells.Item(39, 3).Value = -0.9343346199714787
$ws.Cells.Item(40, 2).Value = 0.3997025630753042
$ws.Cells.Item(40, 3).Value = -0.5442632709721983
$ws.Cells.Item(41, 2).Value = -0.8439833737238367
$ws.Cells.Item(41, 3).Value = -0.6695234394738477
$ws.Cells.Item(42, 2).Value = 1.420299377282191
$ws.Cells.Item(42, 3).Value = -0.05948106311461374
$ws.Cells.Item(43, 2).Value = -1.047372510070521
$ws.Cells.Item(43, 3).Value = 0.4643097799862644
$ws.Cells.Item(44, 2).Value = 0.4492618213325623
$ws.Cells.Item(44, 3).Value = 1.164543727404403
$ws.Cells.Item(45, 2).Value = -1.438805797647627
$ws.Cells.Item(45, 3).Value = 1.118274495760478
$ws.Cells.Item(46, 2).Value = 0.04325537894974833
$ws.Cells.Item(46, 3).Value = 0.7692194348056423
$ws.Cells.Item(47, 2).Value = 0.8948883800423968
$ws.Cells.Item(47, 3).Value = 0.9529118303795087
$ws.Cells.Item(48, 2).Value = 0.7121364849840244
$ws.Cells.Item(48, 3).Value = -0.3280813682164516
$ws.Cells.Item(49, 2).Value = -0.3144437777077964
$ws.Cells.Item(49, 3).Value = -1.037981949723861
$ws.Cells.Item(50, 2).Value = 0.2958377521785464
$ws.Cells.Item(50, 3).Value = 0.6952632191266608
$ws.Cells.Item(51, 2).Value = -0.4683944943133739
$ws.Cells.Item(51, 3).Value = 1.973283416132937
$ws.Cells.Item(52, 2).Value = 0.9765061533991062
$ws.Cells.Item(52, 3).Value = -0.09018112569116743
$ws.Cells.Item(53, 2).Value = -0.8189111800306275
$ws.Cells.Item(53, 3).Value = -0.5348616149443042
$ws.Cells.Item(54, 2).Value = 0.7984055820872634
$ws.Cells.Item(54, 3).Value = 0.7821310160777513
$ws.Cells.Item(55, 2).Value = -0.983976102206072
$ws.Cells.Item(55, 3).Value = -0.144608738519537
$ws.Cells.Item(56, 2).Value = 1.270571652332778
$ws.Cells.Item(56, 3).Value = 1.765715213971919
$ws.Cells.Item(57, 2).Value = -0.006603034504613346
$ws.Cells.Item(57, 3).Value = -0.286719414443142
$ws.Cells.Item(58, 2).Value = -0.4017324373319931
$ws.Cells.Item(58, 3).Value = 0.5249819995714261
$ws.Cells.Item(59, 2).Value = 0.04851654211596454
$ws.Cells.Item(59, 3).Value = -1.512763887543003
$ws.Cells.Item(60, 2).Value = -1.339533302091071
$ws.Cells.Item(60, 3).Value = 0.6738325834995891
$ws.Cells.Item(61, 2).Value = 0.9246830352672981
$ws.Cells.Item(61, 3).Value = -0.1077531474302199
$ws.Cells.Item(62, 2).Value = -0.362427718590265
$ws.Cells.Item(62, 3).Value = 0.2632869405297432
$ws.Cells.Item(63, 2).Value = 0.1186266454692196
$ws.Cells.Item(63, 3).Value = -0.4215271135205001
$ws.Cells.Item(64, 2).Value = -1.188002100962399
$ws.Cells.Item(64, 3).Value = 0.6673674908628129
$ws.Cells.Item(65, 2).Value = -1.880628233888685
$ws.Cells.Item(65, 3).Value = -0.1978679876255
$ws.Cells.Item(66, 2).Value = -1.77824547667701
$ws.Cells.Item(66, 3).Value = -0.9849916834864952
$ws.Cells.Item(67, 2).Value = -0.2613510960543449
$ws.Cells.Item(67, 3).Value = -1.418749517929991
$ws.Cells.Item(68, 2).Value = -0.8594230303627027
$ws.Cells.Item(68, 3).Value = -0.9594827179566222
$ws.Cells.Item(69, 2).Value = -1.502826027541549
$ws.Cells.Item(69, 3).Value = -0.9131374486911226
$ws.Cells.Item(70, 2).Value = -0.8280475372669485
$ws.Cells.Item(70, 3).Value = -0.09310227790034203
$ws.Cells.Item(71, 2).Value = -0.4875854510153331
$ws.Cells.Item(71, 3).Value = -0.7386467929429651
$ws.Cells.Item(72, 2).Value = 2.249225385225403
$ws.Cells.Item(72, 3).Value = 0.6157045521350214
$ws.Cells.Item(73, 2).Value = -0.536050976682898
$ws.Cells.Item(73, 3).Value = 0.8398493684103168
$ws.Cells.Item(74, 2).Value = -0.000259126682991129
$ws.Cells.Item(74, 3).Value = -0.4070462846514331
$ws.Cells.Item(75, 2).Value = -1.345807542313552
$ws.Cells.Item(75, 3).Value = -0.3326405680186354
$ws.Cells.Item(76, 2).Value = -0.8882895590503391
$ws.Cells.Item(76, 3).Value = 0.6005608808997036
$ws.Cells.Item(77, 2).Value = 0.06485664618251867
$ws.Cells.Item(77, 3).Value = -1.848934592689371
$ws.Cells.Item(78, 2).Value = -0.8578388197711537
$ws.Cells.Item(78, 3).Value = -0.3237962582991278
$ws.Cells.Item(79, 2).Value = 0.5109079369671307
$ws.Cells.Item(79, 3).Value = 1.320044902279415
$ws.Cells.Item(80, 2).Value = -0.07868143035309444
$ws.Cells.Item(80, 3).Value = -0.9024747160414333
$ws.Cells.Item(81, 2).Value = 0.7042676453403998
$ws.Cells.Item(81, 3).Value = 0.37959750313202
$ws.Cells.Item(82, 2).Value = 0.1676812134989119
$ws.Cells.Item(82, 3).Value = -0.3774515159114227
$ws.Cells.Item(83, 2).Value = 0.08772380717624008
$ws.Cells.Item(83, 3).Value = -0.3230851048635872
$ws.Cells.Item(84, 2).Value = -0.4055975553642957
$ws.Cells.Item(84, 3).Value = 1.183736920463053
$ws.Cells.Item(85, 2).Value = -0.2536159307640937
$ws.Cells.Item(85, 3).Value = 0.1954726169814617
$ws.Cells.Item(86, 2).Value = -0.6021080182295203
$ws.Cells.Item(86, 3).Value = 0.3503202793467673
$ws.Cells.Item(87, 2).Value = -0.8334508237560108
$ws.Cells.Item(87, 3).Value = -1.152132544257979
$ws.Cells.Item(88, 2).Value = -1.304848672330186
$ws.Cells.Item(88, 3).Value = 0.07460410023412248
$ws.Cells.Item(89, 2).Value = -1.537519732787598
$ws.Cells.Item(89, 3).Value = -0.07064307351042641
$ws.Cells.Item(90, 2).Value = -0.5991997288862755
$ws.Cells.Item(90, 3).Value = 0.1930999299518052
$ws.Cells.Item(91, 2).Value = 1.881840816077361
$ws.Cells.Item(91, 3).Value = -0.1616931509889499
$ws.Cells.Item(92, 2).Value = 0.850594307374502
$ws.Cells.Item(92, 3).Value = -0.130687638173106
$ws.Cells.Item(93, 2).Value = 0.06328323988962521
$ws.Cells.Item(93, 3).Value = -0.2647791550849387
$ws.Cells.Item(94, 2).Value = -0.1134089336824713
$ws.Cells.Item(94, 3).Value = 0.02299509973896885
$ws.Cells.Item(95, 2).Value = 0.9799009882876903
$ws.Cells.Item(95, 3).Value = 0.4164063954553908
$ws.Cells.Item(96, 2).Value = -1.715323537310211
$ws.Cells.Item(96, 3).Value = 1.405001459122823
$ws.Cells.Item(97, 2).Value = 0.1566327268391177
$ws.Cells.Item(97, 3).Value = 0.2459133845969255
$ws.Cells.Item(98, 2).Value = -0.432329694715788
$ws.Cells.Item(98, 3).Value = -0.3299313555613803
$ws.Cells.Item(99, 2).Value = 1.285364363937007
$ws.Cells.Item(99, 3).Value = 0.9261455694758151
$ws.Cells.Item(100, 2).Value = -0.5035318536672438
$ws.Cells.Item(100, 3).Value = 0.3242433330753783
$ws.Cells.Item(101, 2).Value = -0.6592318622190665
$ws.Cells.Item(101, 3).Value = -0.6887154276976315
